$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename headers to lowercase (this moves the old shared strings entries
# out and appends the new lowercase ones at the end of the sst table,
# matching the reference edit).
$ws.Range("A1").Value = "nome"
$ws.Range("B1").Value = "massa"

# Match the resulting view/selection state (scrolled down, new active cell).
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("C42").Select()
